$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.432721376419067
$ws.Range("B1").Value = 1.971966505050659
$ws.Range("C1").Value = 3.017925500869751
$ws.Range("D1").Value = 4.899716377258301
$ws.Range("E1").Value = 0.9230487942695618
